# Shrink the title placeholder on the slide master / "Titel en object" layout
# and let the body / picture placeholders grow upward to reclaim the space
# ("Updated powerpoint template to use less space for the title").
#
# NOTE on the magic point values below: PowerPoint COM works in points
# (1 pt = 12700 EMU) and this host's Left/Top/Width/Height pipeline stores the
# value as a 32-bit float internally, so a naive `emu / 12700.0` can truncate
# one EMU short of the OOXML target after the round-trip back to EMU. Each
# constant here was verified against the produced XML so it lands on the
# exact target EMU.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# --- Slide master: "Title Placeholder 1" ---------------------------------
# <a:off x="838200" y="365126"/> unchanged
# <a:ext cx="10515600" cy="985029"/> -> cy="773863"
$title = $m.Shapes.Item(1)
$title.Height = 60.93410448818898

# --- Slide master: "Text Placeholder 2" (body) ----------------------------
# <a:off x="838200" y="1529542"/> -> y="1267326"
# <a:ext cx="10515600" cy="4647421"/> -> cy="4909637"
$body = $m.Shapes.Item(2)
$body.Top = 99.78944881889764
$body.Height = 386.5855905511811

# --- Layout "1_Titel en object" (CustomLayouts #13): "Picture Placeholder 2"
# <a:off x="838200" y="1803862"/> -> y="1251283"
# <a:ext cx="10517188" cy="4480560"/> -> cy="5470191"
$layout = $m.CustomLayouts.Item(13)
$pic = $layout.Shapes.Item(5)
$pic.Top = 98.52622047244094
$pic.Height = 430.72371078740156
